$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Row 2 originally holds: A2="801193487", B2="11/12/2009", C2="11/12/2009",
    # D2="STATE OF S TEXAS 3H1". B2/C2 stay the same; clear A2 and fix the
    # OCR'd state text in D2.
    $ws.Range("A2").ClearContents()
    $ws.Range("D2").Value = "STATE OF EXAS 3HL"
}
